# Atualiza instrução de trabalho
# - Marca alguns incidentes da planilha ITI como "Resolvido"
# - Adiciona novas linhas de incidentes (semana 32/2025) nas planilhas SPN e ITI

$wb = $excel.ActiveWorkbook

# Helper function to write a full data row (columns A..K) into a worksheet,
# keeping text-like columns (dates, etc.) as plain text instead of letting
# Excel auto-convert them into date serial values.
# NOTE: parameters are passed positionally (named/splatted parameters are
# not reliably bound by this host), so the order below must be respected:
#   ws, row, setor, responsavel, ano, semana, inicioSemana, finalSemana,
#   incidente, backlog, data, status, coordenador
function Add-DataRow {
    param(
        $ws,
        [int]$row,
        [string]$setor,
        [string]$responsavel,
        [int]$ano,
        [int]$semana,
        [string]$inicioSemana,
        [string]$finalSemana,
        [int]$incidente,
        [string]$backlog,
        [string]$data,
        [string]$status,
        [string]$coordenador
    )

    $ws.Cells.Item($row, 1).Value = $setor
    $ws.Cells.Item($row, 2).Value = $responsavel
    $ws.Cells.Item($row, 3).Value = $ano
    $ws.Cells.Item($row, 4).Value = $semana

    $ws.Cells.Item($row, 5).Value = "'" + $inicioSemana
    $ws.Cells.Item($row, 5).Style = "Normal"

    $ws.Cells.Item($row, 6).Value = "'" + $finalSemana
    $ws.Cells.Item($row, 6).Style = "Normal"

    $ws.Cells.Item($row, 7).Value = $incidente

    $ws.Cells.Item($row, 8).Value = "'" + $backlog
    $ws.Cells.Item($row, 8).Style = "Normal"

    $ws.Cells.Item($row, 9).Value = "'" + $data
    $ws.Cells.Item($row, 9).Style = "Normal"

    $ws.Cells.Item($row, 10).Value = $status
    $ws.Cells.Item($row, 11).Value = $coordenador
}

# ---------------------------------------------------------------------------
# Sheet "SPN": add new rows 148-153
# ---------------------------------------------------------------------------
$wsSPN = $wb.Worksheets.Item("SPN")

Add-DataRow $wsSPN 148 "SPN" "Arthur Hassuma" 2025 32 "11/08/2025" "15/08/2025" 342579 "08/2025" "11/08/2025" "Pendente" "Willian Rios"
Add-DataRow $wsSPN 149 "SPN" "Arthur Hassuma" 2025 32 "11/08/2025" "15/08/2025" 343271 "08/2025" "11/08/2025" "Pendente" "Willian Rios"
Add-DataRow $wsSPN 150 "SPN" "Fabio da Silva" 2025 32 "11/08/2025" "15/08/2025" 342380 "08/2025" "11/08/2025" "Pendente" "Willian Rios"
Add-DataRow $wsSPN 151 "SPN" "Higor Jesus"    2025 32 "11/08/2025" "15/08/2025" 342846 "08/2025" "11/08/2025" "Pendente" "Willian Rios"
Add-DataRow $wsSPN 152 "SPN" "Higor Jesus"    2025 32 "11/08/2025" "15/08/2025" 343420 "08/2025" "11/08/2025" "Pendente" "Willian Rios"
Add-DataRow $wsSPN 153 "SPN" "Luana Giese"    2025 32 "11/08/2025" "15/08/2025" 340361 "08/2025" "11/08/2025" "Pendente" "Willian Rios"

# ---------------------------------------------------------------------------
# Sheet "ITI": flip some existing statuses to "Resolvido", add new rows 418-432
# ---------------------------------------------------------------------------
$wsITI = $wb.Worksheets.Item("ITI")

$wsITI.Cells.Item(405, 10).Value = "Resolvido"
$wsITI.Cells.Item(406, 10).Value = "Resolvido"
$wsITI.Cells.Item(411, 10).Value = "Resolvido"
$wsITI.Cells.Item(416, 10).Value = "Resolvido"

Add-DataRow $wsITI 418 "ITI" "Alana Neris"    2025 32 "11/08/2025" "15/08/2025" 342710 "08/2025" "11/08/2025" "Pendente" "Emerson Simette"
Add-DataRow $wsITI 419 "ITI" "Eduardo Batisti" 2025 32 "11/08/2025" "15/08/2025" 342739 "08/2025" "11/08/2025" "Pendente" "Emerson Simette"
Add-DataRow $wsITI 420 "ITI" "Erick da Silva"  2025 32 "11/08/2025" "15/08/2025" 343418 "08/2025" "11/08/2025" "Pendente" "Emerson Simette"
Add-DataRow $wsITI 421 "ITI" "Gabriel López"   2025 32 "11/08/2025" "15/08/2025" 343106 "08/2025" "11/08/2025" "Pendente" "Emerson Simette"
Add-DataRow $wsITI 422 "ITI" "Guilherme Worel" 2025 32 "11/08/2025" "15/08/2025" 343266 "08/2025" "11/08/2025" "Pendente" "Emerson Simette"
Add-DataRow $wsITI 423 "ITI" "Guilherme Worel" 2025 32 "11/08/2025" "15/08/2025" 342874 "08/2025" "11/08/2025" "Pendente" "Emerson Simette"
Add-DataRow $wsITI 424 "ITI" "Guilherme Worel" 2025 32 "11/08/2025" "15/08/2025" 343315 "08/2025" "11/08/2025" "Pendente" "Emerson Simette"
Add-DataRow $wsITI 425 "ITI" "Jorgenaldo Reis" 2025 32 "11/08/2025" "15/08/2025" 343304 "08/2025" "11/08/2025" "Pendente" "Emerson Simette"
Add-DataRow $wsITI 426 "ITI" "Lourival Moizés" 2025 32 "11/08/2025" "15/08/2025" 343006 "08/2025" "11/08/2025" "Pendente" "Emerson Simette"
Add-DataRow $wsITI 427 "ITI" "Lourival Moizés" 2025 32 "11/08/2025" "15/08/2025" 343239 "08/2025" "11/08/2025" "Pendente" "Emerson Simette"
Add-DataRow $wsITI 428 "ITI" "Maria Eduarda"   2025 32 "11/08/2025" "15/08/2025" 343312 "08/2025" "11/08/2025" "Pendente" "Emerson Simette"
Add-DataRow $wsITI 429 "ITI" "Maria Eduarda"   2025 32 "11/08/2025" "15/08/2025" 343313 "08/2025" "11/08/2025" "Pendente" "Emerson Simette"
Add-DataRow $wsITI 430 "ITI" "Maria Eduarda"   2025 32 "11/08/2025" "15/08/2025" 342831 "08/2025" "11/08/2025" "Pendente" "Emerson Simette"
Add-DataRow $wsITI 431 "ITI" "Maria Eduarda"   2025 32 "11/08/2025" "15/08/2025" 342958 "08/2025" "11/08/2025" "Pendente" "Emerson Simette"
Add-DataRow $wsITI 432 "ITI" "Maria Eduarda"   2025 32 "11/08/2025" "15/08/2025" 342914 "08/2025" "11/08/2025" "Pendente" "Emerson Simette"

Write-Host "Linhas adicionadas: SPN 148-153, ITI 418-432. Status atualizado em ITI linhas 405,406,411,416."
